$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = -0.8873889631228897
$ws.Range("J4").Value = 0.4877282578269154
$ws.Range("K4").Value = 0.5675000282201691
$ws.Range("L4").Value = 2.89667423575272
